# Add a new "ammo box" equipment row (row 24) to the equipment record,
# mirroring the formatting of the preceding row (row 23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 23's styled cells (A and E use style index 3)
# onto the corresponding new cells in row 24.
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E23").Copy()
$ws.Range("E24").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's values.
$ws.Range("A24").Value = "ammo box"
$ws.Range("B24").Value = "from Harbor Freight"
$ws.Range("E24").Value = "Golden Home"

# Move the active selection down, as in the source workbook.
$ws.Range("A25").Select()
